$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete footer/notes rows (61-65) that are no longer part of the clean dataset
$ws.Rows("61:65").Delete()

# Rename header row columns to clean snake_case names
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'
$ws.Range("A2").Value = 'Aguascalientes'
$ws.Range("B2").Value = 'Rincón De Romos'
$ws.Range("B3").Value = 'Total'
$ws.Range("A4").Value = 'Baja California'
$ws.Range("B4").Value = 'Mexicali'
$ws.Range("B5").Value = 'Total'
$ws.Range("A6").Value = 'Chihuahua'
$ws.Range("B6").Value = 'Delicias'
$ws.Range("B7").Value = 'Total'
$ws.Range("A8").Value = 'Ciudad De México'
$ws.Range("B8").Value = 'Benito Juárez'
$ws.Range("B9").Value = 'Gustavo A. Madero'
$ws.Range("B10").Value = 'Iztapalapa'
$ws.Range("B11").Value = 'Total'
$ws.Range("A12").Value = 'Coahuila De Zaragoza'
$ws.Range("B12").Value = 'Torreón'
$ws.Range("B13").Value = 'Total'
$ws.Range("A14").Value = 'Colima'
$ws.Range("B14").Value = 'Colima'
$ws.Range("B15").Value = 'Total'
$ws.Range("A16").Value = 'Guerrero'
$ws.Range("B16").Value = 'San Miguel Totolapan'
$ws.Range("B17").Value = 'Tlapa De Comonfort'
$ws.Range("B18").Value = 'Tlapehuala'
$ws.Range("B19").Value = 'Total'
$ws.Range("A20").Value = 'Hidalgo'
$ws.Range("B20").Value = 'Pachuca De Soto'
$ws.Range("B21").Value = 'Total'
$ws.Range("A22").Value = 'Jalisco'
$ws.Range("B22").Value = 'Guadalajara'
$ws.Range("B23").Value = 'Jalostotitlán'
$ws.Range("B24").Value = 'San Miguel El Alto'
$ws.Range("B25").Value = 'Tepatitlán De Morelos'
$ws.Range("B26").Value = 'Zapopan'
$ws.Range("B27").Value = 'Total'
$ws.Range("A28").Value = 'Michoacán De Ocampo'
$ws.Range("B28").Value = 'Acuitzio'
$ws.Range("B29").Value = 'Carácuaro'
$ws.Range("B30").Value = 'Huetamo'
$ws.Range("B31").Value = 'Indaparapeo'
$ws.Range("B32").Value = 'Morelia'
$ws.Range("B33").Value = 'Total'
$ws.Range("A34").Value = 'Nayarit'
$ws.Range("B34").Value = 'Santa María Del Oro'
$ws.Range("B35").Value = 'Total'
$ws.Range("A36").Value = 'Oaxaca'
$ws.Range("B36").Value = 'Santiago Juxtlahuaca'
$ws.Range("B37").Value = 'Total'
$ws.Range("A38").Value = 'Puebla'
$ws.Range("B38").Value = 'Huehuetlán El Grande'
$ws.Range("B39").Value = 'Puebla'
$ws.Range("B40").Value = 'San Pedro Cholula'
$ws.Range("B41").Value = 'Total'
$ws.Range("A42").Value = 'Querétaro'
$ws.Range("B42").Value = 'Tolimán'
$ws.Range("B43").Value = 'Total'
$ws.Range("A44").Value = 'Sinaloa'
$ws.Range("B44").Value = 'Culiacán'
$ws.Range("B45").Value = 'Elota'
$ws.Range("B46").Value = 'Total'
$ws.Range("A47").Value = 'Tamaulipas'
$ws.Range("B47").Value = 'Tampico'
$ws.Range("B48").Value = 'Total'
$ws.Range("A49").Value = 'Veracruz De Ignacio De La Llave'
$ws.Range("B49").Value = 'Coscomatepec'
$ws.Range("B50").Value = 'Ignacio De La Llave'
$ws.Range("B51").Value = 'Sayula De Alemán'
$ws.Range("B52").Value = 'Tlalixcoyan'
$ws.Range("B53").Value = 'Total'
$ws.Range("A54").Value = 'Zacatecas'
$ws.Range("B54").Value = 'Cuauhtémoc'
$ws.Range("B55").Value = 'Guadalupe'
$ws.Range("B56").Value = 'Ojocaliente'
$ws.Range("B57").Value = 'Zacatecas'
$ws.Range("B58").Value = 'Total'
$ws.Range("A59").Value = 'Total'
